# Refactor synthetic array: insert a new "statut_name" column between
# "statut_label" (B) and "NCTId" (C), shifting all subsequent columns
# (NCTId..intervention_type) one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts old C:L (NCTId .. intervention_type)
# to D:M and carries over the bold/bordered header style from the left.
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "statut_name"

# statut_label (column B) -> statut_name (new column C) lookup, derived
# from the four distinct "couleur" statuses used throughout the sheet.
$map = @{
    "noir"   = "pas de résultat ni de publication"
    "rouge"  = "résultat et / ou publication posté"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $map[$label]
}
